$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===========================================================================
# Append 231 new survey responses worth of rows (141-151) below the existing
# 140 data rows, replicating the "Form Responses 1" sheet growth described by
# the commit "more data 231 responses + previous backup".
# ===========================================================================

# --- 1. Formatting ----------------------------------------------------------
# Existing data rows use style s="3" (m/d/yyyy h:mm:ss number format) in
# column A, and style s="4" (general/wrap format) in columns B:K. Copy those
# styles down from row 140 (the current last row) into the new rows. Rows 146
# and 148 have no entry in column G in the source data, so column G format is
# skipped for those two rows, leaving no cell there at all (matching target).
$ws.Range("A140").Copy() | Out-Null
$ws.Range("A141:A151").PasteSpecial(-4122) | Out-Null
$ws.Range("B140").Copy() | Out-Null
$ws.Range("B141:K145").PasteSpecial(-4122) | Out-Null
$ws.Range("B147:K147").PasteSpecial(-4122) | Out-Null
$ws.Range("B149:K151").PasteSpecial(-4122) | Out-Null
$ws.Range("B146:F146").PasteSpecial(-4122) | Out-Null
$ws.Range("H146:K146").PasteSpecial(-4122) | Out-Null
$ws.Range("B148:F148").PasteSpecial(-4122) | Out-Null
$ws.Range("H148:K148").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Cell values -----------------------------------------------------------
# Row 141
$ws.Range("A141").Value = 43759.71853155093
$ws.Range("B141").Value = "Africa"
$ws.Range("C141").Value = 1.0
$ws.Range("D141").Value = "Utility Apps i.e Mobile Banking, Online Shopping, etc"
$ws.Range("E141").Value = 1.0
$ws.Range("F141").Value = 1.0
$ws.Range("G141").Value = "Instagram and mobile banking "
$ws.Range("H141").Value = "Daily"
$ws.Range("I141").Value = "Over 5GB"
$ws.Range("J141").Value = 1.0
$ws.Range("K141").Value = "A day"

# Row 142
$ws.Range("A142").Value = 43759.816405543985
$ws.Range("B142").Value = "Africa"
$ws.Range("C142").Value = 3.0
$ws.Range("D142").Value = "General Browsing"
$ws.Range("E142").Value = 6.0
$ws.Range("F142").Value = 2.0
$ws.Range("G142").Value = "mobile banking-1"
$ws.Range("H142").Value = "Daily"
$ws.Range("I142").Value = "1GB to 3GB"
$ws.Range("J142").Value = 3.0
$ws.Range("K142").Value = "A day"

# Row 143
$ws.Range("A143").Value = 43761.61024681713
$ws.Range("B143").Value = "Africa"
$ws.Range("C143").Value = 1.0
$ws.Range("D143").Value = "General Browsing"
$ws.Range("E143").Value = 4.0
$ws.Range("F143").Value = 2.0
$ws.Range("G143").Value = "1. M-pesa 2. Jumia 3. Uber 4. Mcoop"
$ws.Range("H143").Value = "Monthly"
$ws.Range("I143").Value = "1GB to 3GB"
$ws.Range("J143").Value = 30.0
$ws.Range("K143").Value = "A few Hours"

# Row 144
$ws.Range("A144").Value = 43761.61032630787
$ws.Range("B144").Value = "Africa"
$ws.Range("C144").Value = 1.0
$ws.Range("D144").Value = "General Browsing"
$ws.Range("E144").Value = 4.0
$ws.Range("F144").Value = 2.0
$ws.Range("G144").Value = "1. M-pesa 2. Jumia 3. Uber 4. Mcoop"
$ws.Range("H144").Value = "Monthly"
$ws.Range("I144").Value = "1GB to 3GB"
$ws.Range("J144").Value = 30.0
$ws.Range("K144").Value = "A few Hours"

# Row 145
$ws.Range("A145").Value = 43761.822606435184
$ws.Range("B145").Value = "Africa"
$ws.Range("C145").Value = 5.0
$ws.Range("D145").Value = "Utility Apps i.e Mobile Banking, Online Shopping, etc"
$ws.Range("E145").Value = 10.0
$ws.Range("F145").Value = 3.0
$ws.Range("G145").Value = "Mobile banking aps"
$ws.Range("H145").Value = "Monthly"
$ws.Range("I145").Value = "1GB to 3GB"
$ws.Range("J145").Value = 4.0
$ws.Range("K145").Value = "Several Days"

# Row 146
$ws.Range("A146").Value = 43761.83915774306
$ws.Range("B146").Value = "Africa"
$ws.Range("C146").Value = 2.0
$ws.Range("D146").Value = "General Browsing"
$ws.Range("E146").Value = 4.0
$ws.Range("F146").Value = 1.0
$ws.Range("H146").Value = "Daily"
$ws.Range("I146").Value = "500MB to 1GB"
$ws.Range("J146").Value = 20.0
$ws.Range("K146").Value = "A day"

# Row 147
$ws.Range("A147").Value = 43762.75134616898
$ws.Range("B147").Value = "Middle East"
$ws.Range("C147").Value = 3.0
$ws.Range("D147").Value = "Social Media"
$ws.Range("E147").Value = 5.0
$ws.Range("F147").Value = 2.0
$ws.Range("G147").Value = "Bank, amazon "
$ws.Range("H147").Value = "Monthly"
$ws.Range("I147").Value = "Over 5GB"
$ws.Range("J147").Value = 30.0
$ws.Range("K147").Value = "Automatic reload 0 hours"

# Row 148
$ws.Range("A148").Value = 43762.891612453706
$ws.Range("B148").Value = "Africa"
$ws.Range("C148").Value = 2.0
$ws.Range("D148").Value = "General Browsing"
$ws.Range("E148").Value = 3.0
$ws.Range("F148").Value = 1.0
$ws.Range("H148").Value = "Daily"
$ws.Range("I148").Value = "Less than 100MB"
$ws.Range("J148").Value = 1.0
$ws.Range("K148").Value = "A day"

# Row 149
$ws.Range("A149").Value = 43768.79209307871
$ws.Range("B149").Value = "Africa"
$ws.Range("C149").Value = 2.0
$ws.Range("D149").Value = "Social Media"
$ws.Range("E149").Value = 2.0
$ws.Range("F149").Value = 2.0
$ws.Range("H149").Value = "Monthly"
$ws.Range("I149").Value = "Over 5GB"
$ws.Range("J149").Value = 25.0
$ws.Range("K149").Value = "Several Days"

# Row 150
$ws.Range("A150").Value = 43768.79217355324
$ws.Range("B150").Value = "Africa"
$ws.Range("C150").Value = 2.0
$ws.Range("D150").Value = "Social Media"
$ws.Range("E150").Value = 2.0
$ws.Range("F150").Value = 2.0
$ws.Range("H150").Value = "Monthly"
$ws.Range("I150").Value = "Over 5GB"
$ws.Range("J150").Value = 25.0
$ws.Range("K150").Value = "Several Days"

# Row 151
$ws.Range("A151").Value = 43768.792334328704
$ws.Range("B151").Value = "Africa"
$ws.Range("C151").Value = 2.0
$ws.Range("D151").Value = "Social Media"
$ws.Range("E151").Value = 2.0
$ws.Range("F151").Value = 2.0
$ws.Range("H151").Value = "Monthly"
$ws.Range("I151").Value = "Over 5GB"
$ws.Range("J151").Value = 25.0
$ws.Range("K151").Value = "Several Days"

# --- 3. Multi-line string values ---------------------------------------------
# Assigning a value containing embedded newlines via .Value triggers the
# engine to auto "best fit" the row height (adding an explicit ht/customHeight
# attribute never present in the source file). Writing the text into a scratch
# cell far away, copying it, and pasting *values only* (xlPasteValues) into the
# destination avoids triggering that auto-fit on the destination row; the
# scratch row is removed afterwards so it leaves no trace either.
$ws.Range("Z1000").Value = "Mobile banking`nShopping online"
$ws.Range("Z1001").Value = "Mobile banking`nShopping online"
$ws.Range("Z1002").Value = "Mobile banking`nShopping online"
$ws.Range("Z1000").Copy() | Out-Null
$ws.Range("G149").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1001").Copy() | Out-Null
$ws.Range("G150").PasteSpecial(-4163) | Out-Null
$ws.Range("Z1002").Copy() | Out-Null
$ws.Range("G151").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("Z1000:Z1002").EntireRow.Delete() | Out-Null